$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.297.38'
$ws.Range("E2").Value = '  -1.53%  '
$ws.Range("D3").Value = '2.751.84'
$ws.Range("E3").Value = '  -2.16%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'353.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.22%  '
$ws.Range("D6").Value = "'107.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.32%  '
$ws.Range("D7").Value = "'0.549"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.90%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = "'0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.54%  '
$ws.Range("D10").Value = "'39.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("E11").Value = '  +2.93%  '
$ws.Range("D12").Value = "'0.0834"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.84%  '
$ws.Range("D13").Value = "'19.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").Value = "'7.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.52%  '
$ws.Range("D15").Value = '3.182.27'
$ws.Range("E15").Value = '  -2.26%  '
$ws.Range("D16").Value = '2.750.60'
$ws.Range("E16").Value = '  -2.78%  '
$ws.Range("D17").Value = "'0.921"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").Value = '51.220.06'
$ws.Range("E18").Value = '  -1.31%  '
$ws.Range("D19").Value = "'7.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.47%  '
$ws.Range("E20").Value = '  -2.50%  '
$ws.Range("D21").Value = "'12.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.18%  '
$ws.Range("D22").Value = '0.0₃0959'
$ws.Range("E22").Value = '  -3.18%  '
$ws.Range("D23").Value = "'69.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("D24").Value = "'264.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.46%  '
$ws.Range("D25").Value = "'2.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.53%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = "'25.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.69%  '
$ws.Range("D28").Value = "'0.160"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.48%  '
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").Value = "'10.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("D31").Value = "'6.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.12%  '
$ws.Range("D32").Value = "'34.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.67%  '
$ws.Range("D33").Value = "'51.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.73%  '
$ws.Range("D34").Value = "'0.0441"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.84%  '
$ws.Range("D35").Value = "'0.0830"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("D36").Value = "'5.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.95%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = "'18.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.65%  '
$ws.Range("E39").Value = '  -3.00%  '
$ws.Range("E40").Value = '  -3.97%  '
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("E42").Value = '  -2.08%  '
$ws.Range("E43").Value = '  -2.45%  '
$ws.Range("D44").Value = "'119.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.08%  '
$ws.Range("D45").Value = "'21.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.27%  '
$ws.Range("D46").Value = '2.072.30'
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("D47").Value = "'2.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").Value = "'3.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.96%  '
$ws.Range("E49").Value = '  -3.35%  '
$ws.Range("E50").Value = '  -6.31%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = "'58.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.85%  '
